# "Fruta / hortaliza, semanal" weekly refresh of the price table.
#
# The underlying data source re-ran its "rolling daily logic" and the
# per-date blocks of rows (4-33) got reshuffled: each row's Fecha/Calidad/
# Volumen/Precio mínimo/Precio máximo/Precio promedio ponderado/Precio $/Kg
# (columns D, L, M, N, O, P, S) now hold the values that a *different* row
# held before the refresh - a pure permutation of those seven columns
# across rows 4..33. All other columns (A, B, C, E-K, Q, R, T) are
# untouched.
#
# Snapshot every row's relevant columns first (so the permutation - which
# includes several cycles, not just 2-cycles - doesn't clobber source data
# before it's copied), then write the permuted values back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one record: D, L, M, N, O, P, S
$cols = @(4, 12, 13, 14, 15, 16, 19)

$firstRow = 4
$lastRow = 33

# 1) Snapshot current (pre-edit) values for the columns that change.
$snapshot = @{}
foreach ($r in $firstRow..$lastRow) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) destination row -> source row (the row whose old data it now shows)
$mapping = @{
    4=10; 5=11; 6=17; 7=18; 8=14; 9=15; 10=4; 11=5; 12=6; 13=22; 14=23;
    15=24; 16=25; 17=16; 18=28; 19=29; 20=30; 21=26; 22=27; 23=31; 24=32;
    25=33; 26=21; 27=12; 28=13; 29=19; 30=20; 31=7; 32=8; 33=9
}

# 3) Write each destination row's columns from the snapshotted source row.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcData[$c]
    }
}
